$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B ("sexo") currently holds text values "F"/"M" stored as shared strings.
# Replace them with the numeric codes used in the updated data set: F -> 2, M -> 1.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $txt = $cell.Text
    if ($txt -eq "F") {
        $cell.Value = 2
    } elseif ($txt -eq "M") {
        $cell.Value = 1
    }
}

# Move the active selection, matching the cell that was last selected when the
# workbook was saved.
[void]$ws.Range("G27").Select()
